$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.027338027954102
$ws.Range("B1").Value = 1.256792664527893
$ws.Range("C1").Value = 1.066801905632019
$ws.Range("D1").Value = 1.006998658180237
$ws.Range("E1").Value = 1.060535550117493
